# "Generate Report for Handoff"
# The localization status report is regenerated: the file that was
# "In Translation" (63ace199-...) has now reached "Ready for handoff",
# while cf153297-... remains "In Translation". The report rows are
# re-emitted in (new) status order, so row 2 now carries the
# cf153297 entry and row 3 carries the 63ace199 entry on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "cf153297-95dd-4807-9b07-3dd114d83097.md"
$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"
$ov.Range("D2").Value = "2016-03-22 16:20:10"

$ov.Range("A3").Value = "63ace199-09b8-417d-ab02-ce2988184f6a.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-22 16:21:30"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d35913a084df3bd3290843d2ddf6c864c3ed52d1/e2e/63ace199-09b8-417d-ab02-ce2988184f6a.md", [Type]::Missing, [Type]::Missing, "cf153297-95dd-4807-9b07-3dd114d83097.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d35913a084df3bd3290843d2ddf6c864c3ed52d1/e2e/cf153297-95dd-4807-9b07-3dd114d83097.md", [Type]::Missing, [Type]::Missing, "63ace199-09b8-417d-ab02-ce2988184f6a.md")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "cf153297-95dd-4807-9b07-3dd114d83097.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "In Translation"
$zh.Range("D2").Value = "cf153297-95dd-4807-9b07-3dd114d83097.d7643b751393be2de302d6c844aa042816390d6d.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-22 16:20:01"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = "63ace199-09b8-417d-ab02-ce2988184f6a.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "63ace199-09b8-417d-ab02-ce2988184f6a.27a9d088c84bbe2825ad7efcb1466b1d54a26ebe.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-22 16:21:25"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d35913a084df3bd3290843d2ddf6c864c3ed52d1/e2e/63ace199-09b8-417d-ab02-ce2988184f6a.md", [Type]::Missing, [Type]::Missing, "cf153297-95dd-4807-9b07-3dd114d83097.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ace368b49cb5818aa3fceb0ac818ce554ff4f8ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/63ace199-09b8-417d-ab02-ce2988184f6a.27a9d088c84bbe2825ad7efcb1466b1d54a26ebe.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "cf153297-95dd-4807-9b07-3dd114d83097.d7643b751393be2de302d6c844aa042816390d6d.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d35913a084df3bd3290843d2ddf6c864c3ed52d1/e2e/cf153297-95dd-4807-9b07-3dd114d83097.md", [Type]::Missing, [Type]::Missing, "63ace199-09b8-417d-ab02-ce2988184f6a.md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ace368b49cb5818aa3fceb0ac818ce554ff4f8ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cf153297-95dd-4807-9b07-3dd114d83097.d7643b751393be2de302d6c844aa042816390d6d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "63ace199-09b8-417d-ab02-ce2988184f6a.27a9d088c84bbe2825ad7efcb1466b1d54a26ebe.zh-cn.xlf")

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "cf153297-95dd-4807-9b07-3dd114d83097.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "In Translation"
$de.Range("D2").Value = "cf153297-95dd-4807-9b07-3dd114d83097.d7643b751393be2de302d6c844aa042816390d6d.de-de.xlf"
$de.Range("E2").Value = "2016-03-22 16:20:10"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = "63ace199-09b8-417d-ab02-ce2988184f6a.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "63ace199-09b8-417d-ab02-ce2988184f6a.27a9d088c84bbe2825ad7efcb1466b1d54a26ebe.de-de.xlf"
$de.Range("E3").Value = "2016-03-22 16:21:30"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d35913a084df3bd3290843d2ddf6c864c3ed52d1/e2e/63ace199-09b8-417d-ab02-ce2988184f6a.md", [Type]::Missing, [Type]::Missing, "cf153297-95dd-4807-9b07-3dd114d83097.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cba9fd46e6a972fe5459b98f3a00cc4a628cc0d7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/63ace199-09b8-417d-ab02-ce2988184f6a.27a9d088c84bbe2825ad7efcb1466b1d54a26ebe.de-de.xlf", [Type]::Missing, [Type]::Missing, "cf153297-95dd-4807-9b07-3dd114d83097.d7643b751393be2de302d6c844aa042816390d6d.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d35913a084df3bd3290843d2ddf6c864c3ed52d1/e2e/cf153297-95dd-4807-9b07-3dd114d83097.md", [Type]::Missing, [Type]::Missing, "63ace199-09b8-417d-ab02-ce2988184f6a.md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cba9fd46e6a972fe5459b98f3a00cc4a628cc0d7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cf153297-95dd-4807-9b07-3dd114d83097.d7643b751393be2de302d6c844aa042816390d6d.de-de.xlf", [Type]::Missing, [Type]::Missing, "63ace199-09b8-417d-ab02-ce2988184f6a.27a9d088c84bbe2825ad7efcb1466b1d54a26ebe.de-de.xlf")
